function Set-CellText {
    param($ws, $addr, $text)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws "D2" "61.928.04"
Set-CellText $ws "E2" "  +0.39%  "
Set-CellText $ws "D3" "3.418.65"
Set-CellText $ws "E3" "  +0.48%  "
Set-CellText $ws "E4" "  +0.00%  "
Set-CellText $ws "D5" "409.33"
Set-CellText $ws "E5" "  +1.32%  "
Set-CellText $ws "D6" "128.04"
Set-CellText $ws "E6" "  -2.99%  "
Set-CellText $ws "D7" "0.632"
Set-CellText $ws "E7" "  +7.26%  "
Set-CellText $ws "D8" "1.00"
Set-CellText $ws "E8" "  -0.12%  "
Set-CellText $ws "D9" "0.728"
Set-CellText $ws "E9" "  +6.75%  "
Set-CellText $ws "E10" "  +10.63%  "
Set-CellText $ws "D11" "42.39"
Set-CellText $ws "E11" "  +1.58%  "
Set-CellText $ws "B13" "WrappedliquidstakedEther2.0"
Set-CellText $ws "C13" "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-CellText $ws "D13" "3.960.92"
Set-CellText $ws "E13" "  +0.09%  "
Set-CellText $ws "B14" "Polkadot"
Set-CellText $ws "C14" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText $ws "D14" "9.01"
Set-CellText $ws "E14" "  +8.06%  "
Set-CellText $ws "D15" "21.11"
Set-CellText $ws "E15" "  +7.15%  "
Set-CellText $ws "D16" "0.0000200"
Set-CellText $ws "E16" "  +42.87%  "
Set-CellText $ws "D17" "3.434.01"
Set-CellText $ws "E17" "  +0.26%  "
Set-CellText $ws "D18" "12.26"
Set-CellText $ws "E18" "  +5.54%  "
Set-CellText $ws "D20" "61.849.59"
Set-CellText $ws "E20" "  +0.11%  "
Set-CellText $ws "D21" "441.95"
Set-CellText $ws "E21" "  +42.43%  "
Set-CellText $ws "D22" "91.32"
Set-CellText $ws "E22" "  +10.17%  "
Set-CellText $ws "E23" "  +0.96%  "
Set-CellText $ws "D24" "12.85"
Set-CellText $ws "E24" "  +1.35%  "
Set-CellText $ws "D25" "3.23"
Set-CellText $ws "E25" "  +3.04%  "
Set-CellText $ws "D26" "32.95"
Set-CellText $ws "E26" "  +11.85%  "
Set-CellText $ws "D27" "8.74"
Set-CellText $ws "E27" "  +9.10%  "
Set-CellText $ws "E28" "  -0.76%  "
Set-CellText $ws "B29" "RenderToken"
Set-CellText $ws "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws "D29" "7.54"
Set-CellText $ws "E29" "  -6.54%  "
Set-CellText $ws "B30" "Toncoin"
Set-CellText $ws "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-CellText $ws "D30" "2.72"
Set-CellText $ws "E30" "  -1.06%  "
Set-CellText $ws "D31" "11.91"
Set-CellText $ws "E31" "  +5.21%  "
Set-CellText $ws "E32" "  -0.81%  "
Set-CellText $ws "D33" "0.113"
Set-CellText $ws "E33" "  -0.73%  "
Set-CellText $ws "D34" "42.66"
Set-CellText $ws "E34" "  +0.15%  "
Set-CellText $ws "E35" "  -0.03%  "
Set-CellText $ws "E36" "  +3.22%  "
Set-CellText $ws "D37" "53.21"
Set-CellText $ws "E37" "  +3.62%  "
Set-CellText $ws "E38" "  -0.12%  "
Set-CellText $ws "D39" "3.37"
Set-CellText $ws "E39" "  +0.38%  "
Set-CellText $ws "E40" "  +7.51%  "
Set-CellText $ws "E41" "  -0.61%  "
Set-CellText $ws "D42" "0.313"
Set-CellText $ws "E42" "  -2.80%  "
Set-CellText $ws "D43" "140.96"
Set-CellText $ws "E43" "  +1.05%  "
Set-CellText $ws "D44" "4.22"
Set-CellText $ws "E44" "  +7.35%  "
Set-CellText $ws "D45" "1.97"
Set-CellText $ws "E45" "  +0.70%  "
Set-CellText $ws "E46" "  +8.19%  "
Set-CellText $ws "D47" "16.48"
Set-CellText $ws "E47" "  -0.46%  "
Set-CellText $ws "D48" "22.19"
Set-CellText $ws "E48" "  +4.49%  "
Set-CellText $ws "D49" "3.769.84"
Set-CellText $ws "E49" "  +0.50%  "
Set-CellText $ws "D50" "2.07"
Set-CellText $ws "E50" "  +7.48%  "
Set-CellText $ws "D51" "2.115.50"
Set-CellText $ws "E51" "  +0.73%  "
